# Rename existing sheets
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "getData"
$wb.Worksheets.Item(2).Name = "postData"

# Add new sheet at the end, becomes the active sheet
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "e2eData"

# Column widths (values chosen so the engine's internal char/pixel
# rounding lands on the closest reachable width to the target XML widths:
# 21.28515625, 28.5703125, 17.5703125, 24.7109375)
$ws3.Columns.Item(1).ColumnWidth = 20.5
$ws3.Columns.Item(2).ColumnWidth = 27.666666666666668
$ws3.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws3.Columns.Item(4).ColumnWidth = 23.833333333333336

# Header row
$ws3.Range("A1").Value = "TCs"
$ws3.Range("B1").Value = "uri"
$ws3.Range("C1").Value = "contenttype"
$ws3.Range("D1").Value = "request"

# Row 2 - GET
$ws3.Range("A2").Value = "GET"
$ws3.Range("B2").Value = "/books"
$ws3.Range("C2").Value = "application/json"
$ws3.Range("D2").Value = "NA"

# Row 3 - POST
$postBody = "{`n`n    ""author"": ""Test Name 3"",`n`n    ""isbn"": 103,`n`n    ""name"": ""Test Book 3"",`n`n    ""price"": 999                          }"
$ws3.Range("A3").Value = "POST"
$ws3.Range("B3").Value = "/books"
$ws3.Range("C3").Value = "application/json"
$ws3.Range("D3").Value = $postBody
$ws3.Range("D3").WrapText = $true
$ws3.Rows.Item(3).RowHeight = 105

# Row 4 - PUT
$putBody = "{`n`n    ""author"": ""New Name"",`n`n    ""isbn"": 103,`n`n    ""name"": ""New Name"",`n`n    ""price"": 999                          }"
$ws3.Range("A4").Value = "PUT"
$ws3.Range("B4").Value = "/books"
$ws3.Range("C4").Value = "application/json"
$ws3.Range("D4").Value = $putBody
$ws3.Range("D4").WrapText = $true
$ws3.Rows.Item(4).RowHeight = 90

# Row 5 - DELETE
$ws3.Range("A5").Value = "DELETE"
$ws3.Range("B5").Value = "/books"
$ws3.Range("C5").Value = "application/json"
$ws3.Range("D5").Value = "NA"

# Selection on the new sheet
$ws3.Range("D5").Select()
